$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.931.83"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.640.85"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.25"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5078"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2567"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06384"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.52"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07777"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.298"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "1.644.83"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5458"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "0.0₅7868"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.53"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "25.989.81"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.81"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.430"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.956"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.044"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.881"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.43"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1141"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.881"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.237"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05032"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.264"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.192"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.540"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.363"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8942"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.598"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "1.134.50"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5493"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01556"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("B40").Value = "BabyDogeCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D40").Value = "0.0₈133"
$ws.Range("E40").Value = "  +14.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.542"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.631"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8146"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.92"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "1.780.46"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4529"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.85"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").Value = "  +0.53%  "
